$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row: insert Jan_2026 before Dec_2025, shift Nov_2025 in, drop Oct_2025
$ws.Cells.Item(1, 4).Value = "Jan_2026"
$ws.Cells.Item(1, 5).Value = "Dec_2025"
$ws.Cells.Item(1, 6).Value = "Nov_2025"
$ws.Cells.Item(1, 7).Value = "MoM"
$ws.Cells.Item(1, 8).Value = "QoQ"

# Data rows: rowNum, ISIN, StockName, Jan_2026(D), Dec_2025(E), Nov_2025(F), MoM(G), QoQ(H)
$data = @(
    @(2, 'INE775A01035', 'Samvardhana Motherson International Ltd', 9.708547, 9.34663, 8.757272, 0.361917, 0.951274999999999),
    @(3, 'INE018A01030', 'Larsen & Toubro Limited', 9.707091, 9.13494, 8.796008, 0.5721509999999999, 0.9110829999999996),
    @(4, 'INE814H01029', 'Adani Power Limited', 9.05257, 8.664704, 8.63518, 0.3878659999999989, 0.4173899999999993),
    @(5, 'INE907A01026', 'Kalyani Steels Ltd', 6.182114, 6.077651, 5.868334, 0.104463, 0.3137800000000004),
    @(6, 'INE101I01011', 'Afcons Infrastructure Limited', 5.132306, 5.426595, 5.476188, -0.294289, -0.3438819999999998),
    @(7, 'INE0J1Y01017', 'Life Insurance Corporation Of India', 4.967335, 4.66744, 4.718243, 0.2998950000000002, 0.2490920000000001),
    @(8, 'INE868B01028', 'NCC Ltd', 4.567976, 4.535625, 4.681892, 0.03235100000000024, -0.1139160000000006),
    @(9, 'INE423A01024', 'Adani Enterprises Limited', 4.507944, 4.528577, 1.500787, -0.02063300000000012, 3.007157),
    @(10, 'INE776C01039', 'GMR Airports Limited', 3.37757, 3.515915, 4.741101, -0.1383450000000002, -1.363531),
    @(11, 'INE271C01023', 'DLF Limited', 3.24738, 3.181903, 3.236211, 0.06547700000000001, 0.01116900000000021),
    @(12, 'INE059B01024', 'Simplex Infrastructures Limited', 2.746279, 2.779596, 3.019951, -0.03331700000000026, -0.2736719999999999),
    @(13, 'INE364U01010', 'Adani Green Energy Limited', 2.586681, 2.791836, 1.001686, -0.205155, 1.584995),
    @(14, 'INE347A01017', 'Mangalam Cement Limited', 2.444353, 2.168523, 1.926348, 0.27583, 0.518005),
    @(15, 'INE665A01038', 'Swan Corp Limited', 1.812936, 1.826816, 1.689564, -0.01387999999999989, 0.123372),
    @(16, 'INE464A01036', 'Bharat Bijlee Limited', 1.298112, 1.123699, 1.16802, 0.1744129999999999, 0.1300919999999999),
    @(17, 'INE191B01025', 'Welspun Corp Limited', 1.179162, 1.193344, 1.212454, -0.01418199999999992, -0.03329199999999988),
    @(18, 'INE725E01024', 'The Orissa Minerals Development Co Ltd', 1.110608, 1.124926, 1.065621, -0.01431800000000005, 0.04498700000000011),
    @(19, 'INE239D01028', 'OM INFRA LIMITED', 1.057033, 1.194647, 1.292113, -0.1376140000000001, -0.2350800000000002),
    @(20, 'INE931S01010', 'Adani Energy Solutions Limited', 1.047332, 0.65715, 0.614658, 0.3901819999999999, 0.4326739999999999),
    @(21, 'INE245A01021', 'Tata Power Company Limited', 0.927673, 0.871194, 5.782988, 0.05647899999999995, -4.855314999999999),
    @(22, 'INE095N01031', 'National Building Construction Corp', 0.780386, 0.871083, 0.806517, -0.09069700000000003, -0.02613100000000002),
    @(23, 'INE219X23014', 'India Grid Trust (INVIT)', 0.686119, 0, 0, 0.686119, 0.686119),
    @(24, 'INE206N01018', 'Ravindra Energy Limited', 0.146996, 0.14431, 0.463241, 0.002685999999999994, -0.316245),
    @(25, 'INE467B01029', 'Tata Consultancy Services Limited', 0, 0, 2.951195, 0, -2.951195),
    @(26, 'INE423A20016', 'Adani Enterprises Limited Rights', 0, 0, 0.034942, 0, -0.034942),
    @(27, 'INE002A01018', 'Reliance Industries Limited', 0, 8.524847, 8.221378, -8.524847, -8.221378)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = "quant Infrastructure Fund"
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
}
